# Edit script: add the "macro" named-range/column to the #system sheet of the Nexial
# base-macro library workbook. This mirrors adding MacroLibrary-callable functions
# (description(), expects(var,default), produces(var,value)) as a new "macro" name,
# fixing two pre-existing typos ("runProgram...Parms" and "assertIENavtiveMode"),
# and adding three new Web-automation commands.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1. Insert a new column at M; existing M:AA content shifts right to N:AB
$ws.Columns("M").Insert()

# 2. Header + values for the new "macro" column (M1 header, M2:M4 data)
$ws.Range("M1").Value = "macro"
$colM = @('description()', 'expects(var,default)', 'produces(var,value)')
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $colM[$i]
}

# 3. Column A ("target"): re-sorted list of named-range names, now including "macro"
$colA = @('aws.s3', 'aws.ses', 'base', 'csv', 'desktop', 'excel', 'external', 'image', 'io', 'jms', 'json', 'macro', 'mail', 'number', 'pdf', 'rdbms', 'redis', 'sms', 'sound', 'ssh', 'step', 'web', 'webalert', 'webcookie', 'ws', 'ws.async', 'xml')
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}

# 4. Column H ("external"): fix "runProgram(...Parms)" typo, add "runProgramNoWait(...)"
$colH = @('runJUnit(className)', 'runProgram(programPathAndParams)', 'runProgramNoWait(programPathAndParams)')
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $colH[$i]
}

# 5. Column W ("web", previously V): fix "assertIENavtiveMode()" typo, add "clickOffset(locator,x,y)"
#    and "saveAttributeList(var,locator,attrName)" in their sorted positions
$colW = @('assertAndClick(locator,label)', 'assertAttribute(locator,attrName,value)', 'assertAttributeContains(locator,attrName,contains)', 'assertAttributeNotContains(locator,attrName,contains)', 'assertAttributeNotPresent(locator,attrName)', 'assertAttributePresent(locator,attrName)', 'assertChecked(locator)', 'assertContainCount(locator,text,count)', 'assertCssNotPresent(locator,property)', 'assertCssPresent(locator,property,value)', 'assertElementByAttributes(nameValues)', 'assertElementByText(locator,text)', 'assertElementCount(locator,count)', 'assertElementNotPresent(locator)', 'assertElementPresent(locator)', 'assertFocus(locator)', 'assertFrameCount(count)', 'assertFramePresent(frameName)', 'assertIECompatMode()', 'assertIENativeMode()', 'assertLinkByLabel(label)', 'assertNotChecked(locator)', 'assertNotFocus(locator)', 'assertNotText(locator,text)', 'assertNotVisible(locator)', 'assertOneMatch(locator)', 'assertScrollbarHNotPresent(locator)', 'assertScrollbarHPresent(locator)', 'assertScrollbarVNotPresent(locator)', 'assertScrollbarVPresent(locator)', 'assertTable(locator,row,column,text)', 'assertText(locator,text)', 'assertTextContains(locator,text)', 'assertTextCount(locator,text,count)', 'assertTextList(locator,list,ignoreOrder)', 'assertTextMatches(text,minMatch,scrollTo)', 'assertTextNotPresent(text)', 'assertTextOrder(locator,descending)', 'assertTextPresent(text)', 'assertTitle(text)', 'assertValue(locator,value)', 'assertValueOrder(locator,descending)', 'assertVisible(locator)', 'checkAll(locator)', 'clearLocalStorage()', 'click(locator)', 'clickAndWait(locator,waitMs)', 'clickByLabel(label)', 'clickByLabelAndWait(label,waitMs)', 'clickOffset(locator,x,y)', 'clickWithKeys(locator,keys)', 'close()', 'closeAll()', 'deselect(locator,text)', 'deselectMulti(locator,array)', 'dismissInvalidCert()', 'dismissInvalidCertPopup()', 'doubleClick(locator)', 'doubleClickAndWait(locator,waitMs)', 'doubleClickByLabel(label)', 'doubleClickByLabelAndWait(label,waitMs)', 'dragAndDrop(fromLocator,toLocator)', 'dragTo(fromLocator,xOffset,yOffset)', 'editLocalStorage(key,value)', 'executeScript(var,script)', 'focus(locator)', 'goBack()', 'goBackAndWait()', 'maximizeWindow()', 'mouseOver(locator)', 'open(url)', 'openAndWait(url,waitMs)', 'openHttpBasic(url,username,password)', 'openIgnoreTimeout(url)', 'refresh()', 'refreshAndWait()', 'resizeWindow(width,height)', 'saveAllWindowIds(var)', 'saveAllWindowNames(var)', 'saveAttribute(var,locator,attrName)', 'saveAttributeList(var,locator,attrName)', 'saveCount(var,locator)', 'saveDivsAsCsv(headers,rows,cells,nextPage,file)', 'saveElement(var,locator)', 'saveElements(var,locator)', 'saveLocalStorage(var,key)', 'saveLocation(var)', 'savePageAs(var,sessionIdName,url)', 'savePageAsFile(sessionIdName,url,file)', 'saveTableAsCsv(locator,nextPageLocator,file)', 'saveText(var,locator)', 'saveTextArray(var,locator)', 'saveTextSubstringAfter(var,locator,delim)', 'saveTextSubstringBefore(var,locator,delim)', 'saveTextSubstringBetween(var,locator,start,end)', 'saveValue(var,locator)', 'scrollLeft(locator,pixel)', 'scrollRight(locator,pixel)', 'scrollTo(locator)', 'select(locator,text)', 'selectFrame(locator)', 'selectMulti(locator,array)', 'selectMultiOptions(locator)', 'selectText(locator)', 'selectWindow(winId)', 'selectWindowAndWait(winId,waitMs)', 'selectWindowByIndex(index)', 'selectWindowByIndexAndWait(index,waitMs)', 'toggleSelections(locator)', 'type(locator,value)', 'typeKeys(locator,value)', 'uncheckAll(locator)', 'unselectAllText()', 'upload(fieldLocator,file)', 'verifyContainText(locator,text)', 'verifyText(locator,text)', 'wait(waitMs)', 'waitForElementPresent(locator)', 'waitForPopUp(winId,waitMs)', 'waitForTextPresent(text)', 'waitForTitle(text)')
for ($i = 0; $i -lt $colW.Length; $i++) {
    $ws.Cells.Item($i + 2, 23).Value = $colW[$i]
}

# 6. Update defined-name ranges shifted/extended by the column insert and new rows
$wb.Names.Item('external').RefersTo = '=''#system''!$H$2:$H$4'
$wb.Names.Item('mail').RefersTo = '=''#system''!$N$2:$N$2'
$wb.Names.Item('number').RefersTo = '=''#system''!$O$2:$O$15'
$wb.Names.Item('pdf').RefersTo = '=''#system''!$P$2:$P$16'
$wb.Names.Item('rdbms').RefersTo = '=''#system''!$Q$2:$Q$7'
$wb.Names.Item('redis').RefersTo = '=''#system''!$R$2:$R$10'
$wb.Names.Item('ssh').RefersTo = '=''#system''!$U$2:$U$9'
$wb.Names.Item('step').RefersTo = '=''#system''!$V$2:$V$4'
$wb.Names.Item('target').RefersTo = '=''#system''!$A$2:$A$28'
$wb.Names.Item('web').RefersTo = '=''#system''!$W$2:$W$122'
$wb.Names.Item('webalert').RefersTo = '=''#system''!$X$2:$X$8'
$wb.Names.Item('webcookie').RefersTo = '=''#system''!$Y$2:$Y$8'
$wb.Names.Item('ws').RefersTo = '=''#system''!$Z$2:$Z$17'
$wb.Names.Item('xml').RefersTo = '=''#system''!$AB$2:$AB$13'
$wb.Names.Item('sms').RefersTo = '=''#system''!$S$2:$S$2'
$wb.Names.Item('sound').RefersTo = '=''#system''!$T$2:$T$5'
$wb.Names.Item('ws.async').RefersTo = '=''#system''!$AA$2:$AA$8'

# 7. Add the brand-new "macro" defined name
$wb.Names.Add('macro', '=''#system''!$M$2:$M$4')
